# Updated symbol list with GitHub Actions
# Applies the latest crypto price snapshot to column D (Price) of the
# active worksheet. Values are written as text (using a leading
# apostrophe) so that formatting such as trailing/leading zeros present
# in the source data (e.g. "6.460", "0.00002101") is preserved exactly,
# matching how the values were already stored in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2"  = "245.59"
    "D4"  = "5.362"
    "D5"  = "0.05823"
    "D6"  = "6.463"
    "D7"  = "3.355"
    "D8"  = "0.8109"
    "D9"  = "0.9177"
    "D10" = "0.1402"
    "D11" = "0.07389"
    "D12" = "0.03117"
    "D15" = "3.853"
    "D16" = "0.001555"
    "D17" = "0.04680"
    "D18" = "0.0005972"
    "D19" = "0.006135"
    "D21" = "0.004692"
    "D22" = "0.00008803"
    "D23" = "3.594"
    "D25" = "0.3183"
    "D40" = "0.03849"
    "D41" = "0.006362"
    "D43" = "0.002751"
    "D44" = "0.008543"
    "D45" = "0.00005249"
    "D47" = "0.6862"
    "D48" = "0.001808"
    "D49" = "0.00002101"
    "D50" = "0.0002001"
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = "'" + $updates[$cellRef]
}
